$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shape = $s.Shapes.Item(3)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# 1. bodyPr: normAutofit lnSpcReduction="10000" -> normAutofit (no reduction)
$tf.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

# 2. Split paragraph 2 ("Average Sales from Bottom 3 stores: 5, 33, 44")
#    into two runs: "...5, 33" and ", 44"
$para2 = $tr.Paragraphs(2, 1)
$tail = $para2.Characters(42, 4)   # ", 44"
$tail.Font.Size = 28

# 3. Remove paragraph 3 ("We were not able to correlate the store #'s...")
$para3 = $tr.Paragraphs(3, 1)
$para3.Delete()
